$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2136.5264
$ws.Range("I40").Value = 2136.3635
$ws.Range("J40").Value = 2136.75
$ws.Range("K40").Value = 2136.3635
$ws.Range("L40").Value = 2136.75
$ws.Range("M40").Value = -1961.3635
$ws.Range("N40").Value = -2486.75

$ws.Range("H70").Value = 2233.3333
$ws.Range("I70").Value = 2480
$ws.Range("J70").Value = 1000
$ws.Range("K70").Value = 7440
$ws.Range("L70").Value = 3000
$ws.Range("M70").Value = -7170
$ws.Range("N70").Value = -3540

$ws.Range("H73").Value = 2233.3333
$ws.Range("I73").Value = 2480
$ws.Range("J73").Value = 1000
$ws.Range("K73").Value = 7440
$ws.Range("L73").Value = 3000
$ws.Range("M73").Value = -6504
$ws.Range("N73").Value = -4872

$ws.Range("H112").Value = 927.35
$ws.Range("J112").Value = 939.3158
$ws.Range("L112").Value = 2817.9474
$ws.Range("N112").Value = -5033.9474

$ws.Range("H125").Value = 2565.7144
$ws.Range("I125").Value = 10032
$ws.Range("J125").Value = 1991.3846
$ws.Range("K125").Value = 90288
$ws.Range("L125").Value = 17922.4614
$ws.Range("M125").Value = -87828
$ws.Range("N125").Value = -22842.4614

$ws.Range("H137").Value = 1911.8276
$ws.Range("I137").Value = 1261.091
$ws.Range("K137").Value = 3783.273
$ws.Range("M137").Value = -1233.273

$ws.Range("H138").Value = 7619.1353
$ws.Range("I138").Value = 1372.1875
$ws.Range("J138").Value = 47599.6
$ws.Range("K138").Value = 4116.5625
$ws.Range("L138").Value = 142798.8
$ws.Range("M138").Value = 1023.4375
$ws.Range("N138").Value = -153078.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24512.027
$ws.Range("I32").Value = 4036.9062
$ws.Range("K32").Value = 4036.9062
$ws.Range("M32").Value = -3749.9062

$ws.Range("H45").Value = 72872.78999999999
$ws.Range("I45").Value = 112026.89
$ws.Range("K45").Value = 112026.89
$ws.Range("M45").Value = -111649.89

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 51866.55
$ws.Range("I20").Value = 64600.938
$ws.Range("J20").Value = 929
$ws.Range("K20").Value = 64600.938
$ws.Range("L20").Value = 929
$ws.Range("M20").Value = -64353.938
$ws.Range("N20").Value = -1423

$ws.Range("H64").Value = 698
$ws.Range("I64").Value = 166.33333
$ws.Range("J64").Value = 963.8333
$ws.Range("K64").Value = 166.33333
$ws.Range("L64").Value = 963.8333
$ws.Range("M64").Value = 58.66667000000001
$ws.Range("N64").Value = -1413.8333

$ws.Range("H67").Value = 698
$ws.Range("I67").Value = 166.33333
$ws.Range("J67").Value = 963.8333
$ws.Range("K67").Value = 166.33333
$ws.Range("L67").Value = 963.8333
$ws.Range("M67").Value = 613.6666700000001
$ws.Range("N67").Value = -2523.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29920.428
$ws.Range("I31").Value = 1026.8
$ws.Range("J31").Value = 45972.445
$ws.Range("K31").Value = 1026.8
$ws.Range("L31").Value = 45972.445
$ws.Range("M31").Value = -731.8
$ws.Range("N31").Value = -46562.445

$ws.Range("H34").Value = 29920.428
$ws.Range("I34").Value = 1026.8
$ws.Range("J34").Value = 45972.445
$ws.Range("K34").Value = 1026.8
$ws.Range("L34").Value = 45972.445
$ws.Range("M34").Value = -824.8
$ws.Range("N34").Value = -46376.445

$ws.Range("H86").Value = 2846.7273
$ws.Range("I86").Value = 2599.6667
$ws.Range("J86").Value = 2939.375
$ws.Range("K86").Value = 2599.6667
$ws.Range("L86").Value = 2939.375
$ws.Range("M86").Value = -1476.6667
$ws.Range("N86").Value = -5185.375

$ws.Range("H89").Value = 2846.7273
$ws.Range("I89").Value = 2599.6667
$ws.Range("J89").Value = 2939.375
$ws.Range("K89").Value = 12998.3335
$ws.Range("L89").Value = 14696.875
$ws.Range("M89").Value = -7382.333500000001
$ws.Range("N89").Value = -25928.875

$ws.Range("H132").Value = 28849022
$ws.Range("I132").Value = 25643776
$ws.Range("J132").Value = 38464760
$ws.Range("K132").Value = 76931328
$ws.Range("L132").Value = 115394280
$ws.Range("M132").Value = -76928798
$ws.Range("N132").Value = -115399340

$ws.Range("H134").Value = 806.6
$ws.Range("I134").Value = 713.2353000000001
$ws.Range("J134").Value = 1335.6666
$ws.Range("K134").Value = 2139.7059
$ws.Range("L134").Value = 4006.9998
$ws.Range("M134").Value = 395.2941000000001
$ws.Range("N134").Value = -9076.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 32.166668
$ws.Range("I12").Value = 13.4
$ws.Range("K12").Value = 40.2
$ws.Range("M12").Value = 132.8

$ws.Range("H34").Value = 606.1
$ws.Range("I34").Value = 90.25
$ws.Range("J34").Value = 950
$ws.Range("K34").Value = 270.75
$ws.Range("L34").Value = 2850
$ws.Range("M34").Value = -186.75
$ws.Range("N34").Value = -3018

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 11653.75
$ws.Range("I52").Value = 6030
$ws.Range("J52").Value = 12457.143
$ws.Range("K52").Value = 6030
$ws.Range("L52").Value = 12457.143
$ws.Range("M52").Value = -5771
$ws.Range("N52").Value = -12975.143

$ws.Range("H70").Value = 104390
$ws.Range("I70").Value = 185591
$ws.Range("J70").Value = 5144.3335
$ws.Range("K70").Value = 185591
$ws.Range("L70").Value = 5144.3335
$ws.Range("M70").Value = -185321
$ws.Range("N70").Value = -5684.3335

$ws.Range("H73").Value = 104390
$ws.Range("I73").Value = 185591
$ws.Range("J73").Value = 5144.3335
$ws.Range("K73").Value = 185591
$ws.Range("L73").Value = 5144.3335
$ws.Range("M73").Value = -184655
$ws.Range("N73").Value = -7016.3335

$ws.Range("H132").Value = 3429.2258
$ws.Range("I132").Value = 2319.8667
$ws.Range("J132").Value = 4469.25
$ws.Range("K132").Value = 6959.6001
$ws.Range("L132").Value = 13407.75
$ws.Range("M132").Value = -4429.6001
$ws.Range("N132").Value = -18467.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1014791.2
$ws.Range("I46").Value = 6233.3335
$ws.Range("J46").Value = 1447030.2
$ws.Range("K46").Value = 6233.3335
$ws.Range("L46").Value = 1447030.2
$ws.Range("M46").Value = -6045.3335
$ws.Range("N46").Value = -1447406.2

$ws.Range("H133").Value = 45000
$ws.Range("J133").Value = 45000
$ws.Range("L133").Value = 45000
$ws.Range("N133").Value = -50060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 6824.5
$ws.Range("I54").Value = 7070
$ws.Range("J54").Value = 6775.4
$ws.Range("K54").Value = 7070
$ws.Range("L54").Value = 6775.4
$ws.Range("M54").Value = -6550
$ws.Range("N54").Value = -7815.4

$ws.Range("H62").Value = 6254350
$ws.Range("I62").Value = 62500000
$ws.Range("J62").Value = 4833.3335
$ws.Range("K62").Value = 62500000
$ws.Range("L62").Value = 4833.3335
$ws.Range("M62").Value = -62499376
$ws.Range("N62").Value = -6081.3335

$ws.Range("H65").Value = 6254350
$ws.Range("I65").Value = 62500000
$ws.Range("J65").Value = 4833.3335
$ws.Range("K65").Value = 312500000
$ws.Range("L65").Value = 24166.6675
$ws.Range("M65").Value = -312496880
$ws.Range("N65").Value = -30406.6675

$ws.Range("H132").Value = 2988.9395
$ws.Range("I132").Value = 3106.08
$ws.Range("J132").Value = 2622.875
$ws.Range("K132").Value = 9318.24
$ws.Range("L132").Value = 7868.625
$ws.Range("M132").Value = -6788.24
$ws.Range("N132").Value = -12928.625

$ws.Range("H136").Value = 1030.5186
$ws.Range("I136").Value = 637.8946999999999
$ws.Range("J136").Value = 1963
$ws.Range("K136").Value = 1913.6841
$ws.Range("L136").Value = 5889
$ws.Range("M136").Value = 636.3159000000001
$ws.Range("N136").Value = -10989
